$d = $word.ActiveDocument

# The Bibliografia section currently ends with a trailing block that the
# site rebuild dropped: a blank separator paragraph, a "Ver no Jupiter..."
# paragraph, and a "(c) 2020 ..." footer paragraph. Locate those paragraphs
# by their text (robust to any index drift) and remove them, including
# their paragraph marks, leaving the final blank paragraph and the
# page-break paragraph that follow untouched.

$jupiterFind = $d.Content
[void]$jupiterFind.Find.Execute("Ver no Jupiter Salvar em pdf Salvar em docx", `
    $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$jupiterPara = $jupiterFind.Paragraphs(1)

# The blank paragraph immediately preceding "Ver no Jupiter..." is the
# separator that should go away together with it.
$blankPara = $jupiterPara.Previous()

$copyrightFind = $d.Content
[void]$copyrightFind.Find.Execute([char]0x00A9 + " 2020 . Contact: luizeleno@usp.br. Powered by Jekyll and Github pages. Original theme under Creative Commons Attribution", `
    $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$copyrightPara = $copyrightFind.Paragraphs(1)

$deleteRange = $d.Range($blankPara.Range.Start, $copyrightPara.Range.End)
$deleteRange.Delete()
